$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 243.25
$ws.Range("I33").Value2 = 109.818184
$ws.Range("J33").Value2 = 536.8
$ws.Range("K33").Value2 = 109.818184
$ws.Range("L33").Value2 = 536.8
$ws.Range("M33").Value2 = 119.181816
$ws.Range("N33").Value2 = -994.8
$ws.Range("H111").Value2 = 1490
$ws.Range("J111").Value2 = 2000
$ws.Range("L111").Value2 = 6000
$ws.Range("N111").Value2 = -12134
$ws.Range("H132").Value2 = 1566.3636
$ws.Range("I132").Value2 = 1352.0526
$ws.Range("K132").Value2 = 4056.1578
$ws.Range("M132").Value2 = -1526.1578
$ws.Range("H137").Value2 = 2755.3
$ws.Range("I137").Value2 = 1984
$ws.Range("K137").Value2 = 5952
$ws.Range("M137").Value2 = -3402

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 931.4545000000001
$ws.Range("J2").Value2 = 1036.1
$ws.Range("L2").Value2 = 1036.1
$ws.Range("N2").Value2 = -1262.1
$ws.Range("H32").Value2 = 3268.4062
$ws.Range("I32").Value2 = 3061.3462
$ws.Range("K32").Value2 = 3061.3462
$ws.Range("M32").Value2 = -2774.3462
$ws.Range("H45").Value2 = 3327.077
$ws.Range("I45").Value2 = 2330.8333
$ws.Range("J45").Value2 = 4181
$ws.Range("K45").Value2 = 2330.8333
$ws.Range("L45").Value2 = 4181
$ws.Range("M45").Value2 = -1953.8333
$ws.Range("N45").Value2 = -4935
$ws.Range("H61").Value2 = 6277.778
$ws.Range("I61").Value2 = 5100
$ws.Range("K61").Value2 = 5100
$ws.Range("M61").Value2 = -4888
$ws.Range("H74").Value2 = 1799.1082
$ws.Range("I74").Value2 = 1463.4839
$ws.Range("K74").Value2 = 1463.4839
$ws.Range("M74").Value2 = -589.4838999999999
$ws.Range("H77").Value2 = 1799.1082
$ws.Range("I77").Value2 = 1463.4839
$ws.Range("K77").Value2 = 7317.4195
$ws.Range("M77").Value2 = -2949.4195
$ws.Range("H116").Value2 = 931.4545000000001
$ws.Range("J116").Value2 = 1036.1
$ws.Range("L116").Value2 = 1036.1
$ws.Range("N116").Value2 = -5624.1
$ws.Range("H122").Value2 = 1666.6897
$ws.Range("I122").Value2 = 1291.3182
$ws.Range("K122").Value2 = 3873.9546
$ws.Range("M122").Value2 = -1423.9546
$ws.Range("H136").Value2 = 6277.778
$ws.Range("I136").Value2 = 5100
$ws.Range("K136").Value2 = 15300
$ws.Range("M136").Value2 = -12750

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 931.4545000000001
$ws.Range("J3").Value2 = 1036.1
$ws.Range("L3").Value2 = 1036.1
$ws.Range("N3").Value2 = -1264.1
$ws.Range("H81").Value2 = 106000
$ws.Range("J81").Value2 = 106000
$ws.Range("L81").Value2 = 106000
$ws.Range("N81").Value2 = -108122
$ws.Range("H84").Value2 = 106000
$ws.Range("J84").Value2 = 106000
$ws.Range("L84").Value2 = 318000
$ws.Range("N84").Value2 = -328608
$ws.Range("H86").Value2 = 4610.3076
$ws.Range("I86").Value2 = 3228
$ws.Range("K86").Value2 = 3228
$ws.Range("M86").Value2 = -2105
$ws.Range("H89").Value2 = 4610.3076
$ws.Range("I89").Value2 = 3228
$ws.Range("K89").Value2 = 16140
$ws.Range("M89").Value2 = -10524
$ws.Range("H94").Value2 = 1248
$ws.Range("I94").Value2 = 1239.5555
$ws.Range("K94").Value2 = 1239.5555
$ws.Range("M94").Value2 = -788.5554999999999
$ws.Range("H107").Value2 = 4579.923
$ws.Range("I107").Value2 = 3092.875
$ws.Range("J107").Value2 = 6959.2
$ws.Range("K107").Value2 = 3092.875
$ws.Range("L107").Value2 = 6959.2
$ws.Range("M107").Value2 = -1172.875
$ws.Range("N107").Value2 = -10799.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2403.3
$ws.Range("I58").Value2 = 1488
$ws.Range("J58").Value2 = 3318.6
$ws.Range("K58").Value2 = 1488
$ws.Range("L58").Value2 = 3318.6
$ws.Range("M58").Value2 = -1285
$ws.Range("N58").Value2 = -3724.6
$ws.Range("H86").Value2 = 9300
$ws.Range("H89").Value2 = 9300
$ws.Range("H107").Value2 = 665.4
$ws.Range("I107").Value2 = 516.6667
$ws.Range("K107").Value2 = 516.6667
$ws.Range("M107").Value2 = 1403.3333
$ws.Range("H136").Value2 = 2403.3
$ws.Range("I136").Value2 = 1488
$ws.Range("J136").Value2 = 3318.6
$ws.Range("K136").Value2 = 4464
$ws.Range("L136").Value2 = 9955.799999999999
$ws.Range("M136").Value2 = -1914
$ws.Range("N136").Value2 = -15055.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 143.6923
$ws.Range("I12").Value2 = 16.666666
$ws.Range("K12").Value2 = 49.999998
$ws.Range("M12").Value2 = 123.000002
$ws.Range("H38").Value2 = 380.5909
$ws.Range("I38").Value2 = 362.78946
$ws.Range("J38").Value2 = 493.33334
$ws.Range("K38").Value2 = 1088.36838
$ws.Range("L38").Value2 = 1480.00002
$ws.Range("M38").Value2 = -741.3683800000001
$ws.Range("N38").Value2 = -2174.00002
$ws.Range("H107").Value2 = 457
$ws.Range("J107").Value2 = 900
$ws.Range("L107").Value2 = 2700
$ws.Range("N107").Value2 = -6540
$ws.Range("H136").Value2 = 6679.143
$ws.Range("I136").Value2 = 2499.5
$ws.Range("K136").Value2 = 7498.5
$ws.Range("M136").Value2 = -2398.5
$ws.Range("H138").Value2 = 3289
$ws.Range("I138").Value2 = 1712.8572
$ws.Range("K138").Value2 = 5138.571599999999
$ws.Range("M138").Value2 = 1.428400000000693
$ws.Range("H140").Value2 = 1586.5883
$ws.Range("I140").Value2 = 1264.8
$ws.Range("K140").Value2 = 3794.4
$ws.Range("M140").Value2 = 1385.6
$ws.Range("H141").Value2 = 3532.25
$ws.Range("I141").Value2 = 1376.3334
$ws.Range("K141").Value2 = 4129.0002
$ws.Range("M141").Value2 = 1050.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 10212.286
$ws.Range("I80").Value2 = 8298.4
$ws.Range("K80").Value2 = 8298.4
$ws.Range("M80").Value2 = -7300.4
$ws.Range("H83").Value2 = 10212.286
$ws.Range("I83").Value2 = 8298.4
$ws.Range("K83").Value2 = 41492
$ws.Range("M83").Value2 = -36500
$ws.Range("H132").Value2 = 2942.0527
$ws.Range("I132").Value2 = 2344.7693
$ws.Range("K132").Value2 = 7034.3079
$ws.Range("M132").Value2 = -4504.3079

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 5834
$ws.Range("I61").Value2 = 2502
$ws.Range("J61").Value2 = 7500
$ws.Range("K61").Value2 = 2502
$ws.Range("L61").Value2 = 7500
$ws.Range("M61").Value2 = -2300
$ws.Range("N61").Value2 = -7904
$ws.Range("H68").Value2 = 8637.833000000001
$ws.Range("J68").Value2 = 8712.714
$ws.Range("L68").Value2 = 8712.714
$ws.Range("N68").Value2 = -10210.714
$ws.Range("H71").Value2 = 8637.833000000001
$ws.Range("J71").Value2 = 8712.714
$ws.Range("L71").Value2 = 43563.57
$ws.Range("N71").Value2 = -51051.57
$ws.Range("H113").Value2 = 5834
$ws.Range("I113").Value2 = 2502
$ws.Range("J113").Value2 = 7500
$ws.Range("K113").Value2 = 2502
$ws.Range("L113").Value2 = 7500
$ws.Range("M113").Value2 = -332
$ws.Range("N113").Value2 = -11840
$ws.Range("H136").Value2 = 4700
$ws.Range("I136").Value2 = 3266.6667
$ws.Range("K136").Value2 = 9800.000100000001
$ws.Range("M136").Value2 = -7250.000100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value2 = 39998
$ws.Range("J69").Value2 = 39998
$ws.Range("L69").Value2 = 39998
$ws.Range("N69").Value2 = -41496
$ws.Range("H72").Value2 = 39998
$ws.Range("J72").Value2 = 39998
$ws.Range("L72").Value2 = 119994
$ws.Range("N72").Value2 = -127482
$ws.Range("H81").Value2 = 9499.143
$ws.Range("I81").Value2 = 7666.6665
$ws.Range("K81").Value2 = 15333.333
$ws.Range("M81").Value2 = -14272.333
$ws.Range("H84").Value2 = 9499.143
$ws.Range("I84").Value2 = 7666.6665
$ws.Range("K84").Value2 = 76666.66500000001
$ws.Range("M84").Value2 = -71362.66500000001
$ws.Range("H122").Value2 = 4278.1875
$ws.Range("I122").Value2 = 4267.9287
$ws.Range("K122").Value2 = 12803.7861
$ws.Range("M122").Value2 = -10353.7861
$ws.Range("H126").Value2 = 5556.316
$ws.Range("I126").Value2 = 2973
$ws.Range("K126").Value2 = 8919
$ws.Range("M126").Value2 = -6449
$ws.Range("H141").Value2 = 299997.5
$ws.Range("J141").Value2 = 99995
$ws.Range("L141").Value2 = 99995
$ws.Range("N141").Value2 = -110355
